$d = $word.ActiveDocument

function Insert-PlainParagraph($text, $bold) {
    $trailing = $d.Paragraphs.Item($d.Paragraphs.Count)
    $ip = $d.Range($trailing.Range.Start, $trailing.Range.Start)
    $ip.InsertBefore($text + "`r")
    $p = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
    if ($bold) {
        $p.Range.Bold = 1
    }
    return $p
}

function Insert-BulletParagraph($text, $donorIndex) {
    $trailing = $d.Paragraphs.Item($d.Paragraphs.Count)
    $ip = $d.Range($trailing.Range.Start, $trailing.Range.Start)
    $donor = $d.Paragraphs.Item($donorIndex)
    $ip.FormattedText = $donor.Range.FormattedText
    $p = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
    $p.Range.Text = $text
    return $p
}

# 1) Bold heading: "Feedback from Gergana on full plan to consider in future"
Insert-PlainParagraph "Feedback from Gergana on full plan to consider in future" $true | Out-Null

# 2) Blank ListParagraph line (no numbering)
$trailing = $d.Paragraphs.Item($d.Paragraphs.Count)
$ip = $d.Range($trailing.Range.Start, $trailing.Range.Start)
$ip.InsertBefore("`r")
$pBlank = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$pBlank.Style = "List Paragraph"
$pBlank.Range.Font.NameAscii = "Arial"
$pBlank.Range.Font.NameOther = "Arial"
$pBlank.Range.Font.NameBi = "Arial"
$pBlank.Range.LanguageID = "en-US"

# Locate an existing numId=3/ilvl=0 ListParagraph bullet to use as a donor for formatting
$donorIndex = 43  # "Better to fully answer first question..." paragraph

# 3) Bulleted items (numId=3, ilvl=0)
Insert-BulletParagraph "Use active voice" $donorIndex | Out-Null
Insert-BulletParagraph "Make sure first and last sentences are stand-alone and link together" $donorIndex | Out-Null
Insert-BulletParagraph "Can make a conceptual diagram about objectives" $donorIndex | Out-Null
Insert-BulletParagraph "Even if one reference has everything, use a variety of new and old papers" $donorIndex | Out-Null
Insert-BulletParagraph "Potential to make hypotheses even more specific/directional – something to think about" $donorIndex | Out-Null
Insert-BulletParagraph "Split up methods by sub-question to make clearer" $donorIndex | Out-Null
Insert-BulletParagraph "Add model equations" $donorIndex | Out-Null

# Now remove the old trailing blank paragraph (the one that always remains at the very end)
$trailing = $d.Paragraphs.Item($d.Paragraphs.Count)
$trailing.Range.Delete() | Out-Null

Write-Host "Final paragraph count:" $d.Paragraphs.Count
